$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("H6")
$r.Font.Name = "Roboto"
$r.Font.Color = 16711680
$r.Font.Underline = 2
$r.Interior.Color = 3816153
$r.Borders.LineStyle = 1
$r.Borders.Color = 3816153
